$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.608.24"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "2.490.58"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "491.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.94"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +8.78%  "
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").Value = "2.503.55"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  +5.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0983"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "2.923.66"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "56.670.34"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.22"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.91%  "
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "2.500.89"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("E19").Value = "  +4.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.27"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.44"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.90"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.86"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.65%  "
$ws.Range("E25").Value = "  +1.99%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.163"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.995"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Value = "2.597.27"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("E29").Value = "  +4.04%  "
$ws.Range("D30").Value = "0.0₃0801"
$ws.Range("E30").Value = "  +4.15%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.24"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.38"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("E34").Value = "  +3.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.21"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.16"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.01%  "
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("E38").Value = "  +3.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.39"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.92"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.51"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.99%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.613"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0559"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.88"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +9.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "263.52"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0926"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.49%  "
$ws.Range("E48").Value = "  +3.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.21"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.75"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("D51").Value = "1.903.53"
$ws.Range("E51").Value = "  -3.37%  "
